$wb = $excel.ActiveWorkbook

# --- Sheet 2 (tab) rename: "Include from index.htm" -> "Include #0" ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Include #0"

# --- Sheet 1 ("Metadata") edits ---
$ws1 = $wb.Worksheets.Item(1)

# 1) Insert a new "Identifier" row right after the "URL" row (before the old "Version" row, row 3).
$ws1.Rows.Item(3).Insert()
$ws1.Range("A2:B2").Copy()
$ws1.Range("A3:B3").PasteSpecial(-4122)
$ws1.Range("A3").Value = "Identifier"
$ws1.Range("B3").Value = "OID:1.3.6.1.4.1.19376.1.5.3.1.3.43.48.2"

# 2) Version value: "1.0.0-comment" -> "1.0.0" (now row 4 after the insert above).
$ws1.Range("B4").Value = "1.0.0"

# 3) Date value update (now row 9).
$ws1.Range("B9").Value = "2024-12-04T15:50:20-06:00"

# 4) Publisher value update (now row 10).
$ws1.Range("B10").Value = "IHE Patient Care Coordination Committee"

# 5) Contact: the single "Contact"/"No display for ContactDetail" row (now row 11) becomes
#    three separate "Contact" rows with different values. Insert two more rows below it.
$ws1.Range("B11").Value = "null (https://www.ihe.net/ihe_domains/patient_care_coordination/)"

$ws1.Rows.Item(12).Insert()
$ws1.Range("A11:B11").Copy()
$ws1.Range("A12:B12").PasteSpecial(-4122)
$ws1.Range("A12").Value = "Contact"
$ws1.Range("B12").Value = "null (pcc@ihe.net)"

$ws1.Rows.Item(13).Insert()
$ws1.Range("A12:B12").Copy()
$ws1.Range("A13:B13").PasteSpecial(-4122)
$ws1.Range("A13").Value = "Contact"
$ws1.Range("B13").Value = "IHE Patient Care Coordination Committee (pcc@ihe.net)"
